# Insert a new weekly price-report row above row 19 (pushing the existing
# rows 19-133 down to 20-134) and populate it with the new week's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value = 10
$ws.Cells.Item(19, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value = "La Araucanía"
$ws.Cells.Item(19, 4).Value = 44473
$ws.Cells.Item(19, 5).Value = 9
$ws.Cells.Item(19, 6).Value = 100112005
$ws.Cells.Item(19, 7).Value = "Puerro"
$ws.Cells.Item(19, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 50
$ws.Cells.Item(19, 11).Value = 7000
$ws.Cells.Item(19, 12).Value = 7000
$ws.Cells.Item(19, 13).Value = 7000
$ws.Cells.Item(19, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(19, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(19, 16).Value = 583
$ws.Cells.Item(19, 17).Value = 12
$ws.Cells.Item(19, 18).Value = "Hortaliza"
